# Update cryptos price list (Price + Volume(1h) columns, and two swapped rows).
# Note: "Price" column (D) values are stored as text in the workbook, and several
# look like plain numbers (e.g. "1.00", "0.0000237", "0.290"). Assigning them with a
# leading apostrophe forces Excel to keep them as literal text (quotePrefix) instead
# of silently re-interpreting/reformatting them as numbers (which would drop
# significant trailing zeros or switch to scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.703.96"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'3.078.45"
$ws.Range("E3").Value = "  -11.61%  "
$ws.Range("E4").Value = "  -2.73%  "
$ws.Range("D5").Value = "'588.12"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").Value = "'155.27"
$ws.Range("E6").Value = "  +4.47%  "
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "'3.077.82"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("E10").Value = "  -4.59%  "
$ws.Range("D11").Value = "'5.86"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("D14").Value = "'36.72"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "'3.586.25"
$ws.Range("E16").Value = "  -11.41%  "
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").Value = "'63.653.36"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'3.079.44"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").Value = "'468.99"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "'14.28"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").Value = "'0.704"
$ws.Range("E22").Value = "  -5.96%  "
$ws.Range("D23").Value = "'7.49"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").Value = "'12.83"
$ws.Range("E25").Value = "  -5.83%  "
$ws.Range("D26").Value = "'80.41"
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "'7.45"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.66"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("E32").Value = "  -6.24%  "
$ws.Range("E33").Value = "  -8.00%  "
$ws.Range("D34").Value = "'27.07"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").Value = "'0.0₃0827"
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("E39").Value = "  -5.94%  "
$ws.Range("D40").Value = "'50.57"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").Value = "'9.11"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").Value = "'432.42"
$ws.Range("E42").Value = "  -8.34%  "
$ws.Range("D43").Value = "'0.290"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0360"
$ws.Range("E45").Value = "  -4.67%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "'39.87"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'2.811.38"
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").Value = "'129.95"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D50").Value = "'24.86"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "'2.21"
$ws.Range("E51").Value = "  -4.01%  "
